# Insert a new weekly price record as row 160 (pushing the existing rows
# 160-226 down to 161-227), matching the "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 160..226 down by one to make room for the new record.
$ws.Rows.Item(160).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A160").Value = 4
$ws.Range("B160").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C160").Value = "Los Lagos"
$ws.Range("D160").Value = 44726
$ws.Range("E160").Value = 10
$ws.Range("F160").Value = 100112039
$ws.Range("G160").Value = "Ciboulette"
$ws.Range("H160").Value = "Sin especificar"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 240
$ws.Range("K160").Value = 2500
$ws.Range("L160").Value = 2500
$ws.Range("M160").Value = 2500
$ws.Range("N160").Value = "`$/docena de atados"
$ws.Range("O160").Value = "Región Metropolitana"
$ws.Range("P160").Value = 833
$ws.Range("Q160").Value = 3
$ws.Range("R160").Value = "Hortaliza"
